$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows to append (update through 09/09, as per commit message
# "aggiornamento a 9/09 compreso"). Columns: date serial, nuovi pos.,
# somma mobile 7gg., somma mobile 7gg. per 100mila abitanti.
$data = @(
    @(44441, 1, 1, 53.73455131649651),
    @(44442, 0, 1, 53.73455131649651),
    @(44443, 0, 1, 53.73455131649651),
    @(44444, 0, 1, 53.73455131649651),
    @(44445, 0, 1, 53.73455131649651),
    @(44446, 0, 1, 53.73455131649651),
    @(44447, 0, 1, 53.73455131649651),
    @(44448, 0, 0, 0)
)

$lastRow = 366
$row = $lastRow + 1
foreach ($d in $data) {
    # Copy the formatting (date style) of the last existing data row down
    # onto the new row before writing values into it.
    $ws.Range("A$lastRow").Copy($ws.Range("A$row"))

    $ws.Range("A$row").Value = $d[0]
    $ws.Range("B$row").Value = $d[1]
    $ws.Range("C$row").Value = $d[2]
    $ws.Range("D$row").Value = $d[3]

    $row = $row + 1
}
